$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.464.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.48%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.841.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5203"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3276"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06778"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7704"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07690"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.823.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.029"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.38%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007950"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.447.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.061.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.576"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.489"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.967"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.210"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.644"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.196"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.133"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08724"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04805"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.835"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7078"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.222"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01761"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4839"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8958"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.065"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.728"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4164"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05871"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.980"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1219"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8882"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.74%  "
